# Apply the data fixes described by the commit "fix constraint time 2 3"
# across the logistic routing workbook.

$wb = $excel.ActiveWorkbook

# --- Time Matrix: make travel time between Depot and No.1 symmetric (8 -> 3) ---
$wsTime = $wb.Worksheets.Item("Time Matrix")
$wsTime.Activate()
$wsTime.Range("C2").Value = 3
$wsTime.Range("B3").Value = 3
$wsTime.Range("B3").Select()

# --- Time Window: update the earliest/latest service start times ---
$wsWindow = $wb.Worksheets.Item("Time Window")
$wsWindow.Activate()
$wsWindow.Range("C2").Value = 3
$wsWindow.Range("B3").Value = 10
$wsWindow.Range("B4").Value = 16
$wsWindow.Range("C4").Value = 24
$wsWindow.Range("B5").Value = 12
$wsWindow.Range("B7").Value = 12
$wsWindow.Range("C8").Value = 12
$wsWindow.Range("B10").Value = 11
$wsWindow.Range("C10").Value = 16
$wsWindow.Range("C11").Value = 9
$wsWindow.Range("B3").Select()

# --- Demand Matrix: update customer demand values ---
$wsDemand = $wb.Worksheets.Item("Demand Matrix")
$wsDemand.Activate()
$wsDemand.Range("B2").Value = 30
$wsDemand.Range("B3").Value = 42
$wsDemand.Range("B6").Value = 68
$wsDemand.Range("B7").Value = 35
$wsDemand.Range("B3").Select()

# --- Capicity: update vehicle capacities (this ends up the active tab) ---
$wsCap = $wb.Worksheets.Item("Capicity")
$wsCap.Activate()
$wsCap.Range("B2").Value = 67
$wsCap.Range("B3").Value = 50
$wsCap.Range("B4").Value = 45
$wsCap.Range("B5").Value = 150
$wsCap.Range("B2").Select()

# --- Time Service: just move the selection, no data changes ---
$wsService = $wb.Worksheets.Item("Time Service")
$wsService.Activate()
$wsService.Range("B11").Select()

# Leave "Capicity" as the final active sheet, matching the saved workbook state.
$wsCap.Activate()
